# [PV-94][WIP] Support for plans without sticky-ids or levels
# Rename the plan-import column headers so the sheet no longer assumes
# "sticky ids" / "Name" / "Start" / "Finish" terminology.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-01")

$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Match the author's resulting selection state.
$ws.Range("F1").Select()
